$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1074
$ws1.Range("F5").Value = 2932
$ws1.Range("F12").Value = 159
$ws1.Range("F13").Value = 73
$ws1.Range("F14").Value = 2770
$ws1.Range("F15").Value = 1047

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1074
$ws4.Range("F6").Value = 2932
$ws4.Range("F14").Value = 159
$ws4.Range("F15").Value = 73
$ws4.Range("F16").Value = 2770
$ws4.Range("F17").Value = 1047
